$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.105.22'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.656.54'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '218.37'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = '0.5304'
$ws.Range('E6').Value = '  +1.59%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.2613'
$ws.Range('E8').Value = '  -2.24%  '
$ws.Range('D9').Value = '0.06340'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '20.39'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D11').Value = '0.07762'
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = '4.499'
$ws.Range('E12').Value = '  +1.57%  '
$ws.Range('D13').Value = '1.655.80'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '0.5471'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '0.0₅8151'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '65.29'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').Value = '26.133.05'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').Value = '4.540'
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('D20').Value = '193.77'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '10.08'
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').Value = '6.007'
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '140.51'
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('D25').Value = '0.1242'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = '7.281'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').Value = '16.15'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').Value = '0.05950'
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '1.278'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').Value = '3.511'
$ws.Range('E31').Value = '  -3.62%  '
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('D33').Value = '1.544'
$ws.Range('E33').Value = '  -5.34%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '0.9466'
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.413'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').Value = '2.762'
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('D37').Value = '0.5639'
$ws.Range('E37').Value = '  -4.37%  '
$ws.Range('D38').Value = '0.01609'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').Value = '5.849'
$ws.Range('E39').Value = '  -1.68%  '
$ws.Range('D40').Value = '0.8473'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.013.50'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '101.27'
$ws.Range('E43').Value = '  +1.52%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = '56.96'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('E46').Value = '  -6.24%  '
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').Value = '0.4287'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = '0.05153'
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '7.728'
$ws.Range('E51').Value = '  -4.45%  '
